$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 12.944643
$ws.Range("N2").Value = 38.833929
$ws.Range("O2").Value = 0.3748900893017936
$ws.Range("P2").Value = 0.3748900893017936
$ws.Range("Q2").Value = 0.602819079867
$ws.Range("R2").Value = 5.425371718803
$ws.Range("S2").Value = 0.3748900893017936
$ws.Range("T2").Value = 0.3748900893017936

# Row 3
$ws.Range("O3").Value = 0.4353841324781401
$ws.Range("P3").Value = 0.4353841324781401
$ws.Range("S3").Value = 0.4353841324781401
$ws.Range("T3").Value = 0.4353841324781401

# Row 4
$ws.Range("O4").Value = 0.1897257782200662
$ws.Range("P4").Value = 0.1897257782200662
$ws.Range("S4").Value = 0.1897257782200662
$ws.Range("T4").Value = 0.1897257782200662
